# FreeCRMTestDataSheet.xlsx - "Committing minor changes to fix test NG errors"
# Update the Contacts sheet: swap in new First/Last names for the existing
# three contacts and add a brand-new fourth contact row.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Contacts")
$ws.Activate()

# Row 2: Dr. Matt Williams -> Dr. Gunit Thapar (same company)
$ws.Range("B2").Value = "Gunit"
$ws.Range("C2").Value = "Thapar"

# Row 3: Mr. Mark Keller -> Mr. Asheen Antony (same company)
$ws.Range("B3").Value = "Asheen"
$ws.Range("C3").Value = "Antony"

# Row 4: Mrs. Anshika Khandelwal -> Mrs. Manpreet Singh (same company)
$ws.Range("B4").Value = "Manpreet"
$ws.Range("C4").Value = "Singh"

# Row 5: brand-new contact
$ws.Range("A5").Value = "Dr."
$ws.Range("B5").Value = "Rahul"
$ws.Range("C5").Value = "Varma"
$ws.Range("D5").Value = "Limestone Pvt. Ltd."

# Update the active selection on the Contacts sheet to B2
[void]$ws.Range("B2").Select()
